$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of connector words ("de"->"De", "el"->"El", "del"->"Del", "los"->"Los")
$ws.Range("B10").Value = "Comitán De Domínguez"
$ws.Range("B13").Value = "Hidalgo Del Parral"
$ws.Range("A19").Value = "Ciudad De México"
$ws.Range("A32").Value = "Estado De México"
$ws.Range("B37").Value = "Apaseo El Alto"
$ws.Range("B47").Value = "Coyuca De Benítez"
$ws.Range("B48").Value = "Cuetzala Del Progreso"
$ws.Range("B55").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B58").Value = "Tlajomulco De Zúñiga"
$ws.Range("B66").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B91").Value = "San Juan Del Río"

# Remove the trailing footer/metadata rows (117-121), leaving row 115 as the last row
$ws.Rows.Item(117).Resize(5).Delete()
